$d = $word.ActiveDocument

# Locate the "11-Oct" date text (due date of the "Preliminary Architecture" row).
# Find.Execute mutates the Range it is called on to the bounds of the match.
$hit = $d.Range(0, $d.Content.End)
$found = $hit.Find.Execute("11-Oct", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Narrow to just the "11" part (the first two characters of the match)
    # and change it to "12".
    $num = $d.Range($hit.Start, $hit.Start + 2)
    $num.Text = "12"

    # Re-acquire the (now 2-char) range for "12" and force Word to materialize
    # it as its own run (distinct from the following "-Oct" text) by toggling
    # a character-format property on it and back off again.
    $num2 = $d.Range($hit.Start, $hit.Start + 2)
    $num2.Bold = 1
    $num2.Bold = 0
}
